# Mark the last functions (ft_putchar_fd, ft_putstr_fd, ft_putendl_fd,
# ft_putnbr_fd) as Done/Valid to match mandatory part 1 and 2, as per
# commit message "add last functions to match mandatory part 1 and 2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E51:E54").Value = "Done"
$ws.Range("F51:F54").Value = "Valid"
